$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the order of "Venezuela" and "Estado de Palestina" rows ---
# Before: row 124 = Venezuela, row 125 = Estado de Palestina
# After:  row 124 = Estado de Palestina (with new data), row 125 = Venezuela (old Venezuela data)
$ws.Range("A124").Value = "Estado de Palestina"
$ws.Range("A125").Value = "Venezuela"

# --- Update statistic values per the diff ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1202911
$ws.Range("C4").Value = 14789
$ws.Range("D4").Value = 181901
$ws.Range("E4").Value = 951842
$ws.Range("G4").Value = 570
$ws.Range("H4").Value = 69168

# Row 17: Peru
$ws.Range("B17").Value = 47372
$ws.Range("C17").Value = 1444
$ws.Range("D17").Value = 14427
$ws.Range("E17").Value = 31601
$ws.Range("F17").Value = 694
$ws.Range("G17").Value = 58
$ws.Range("H17").Value = 1344

# Row 18: India
$ws.Range("B18").Value = 46437
$ws.Range("C18").Value = 3932
$ws.Range("D18").Value = 12842
$ws.Range("E18").Value = 32029
$ws.Range("G18").Value = 175
$ws.Range("H18").Value = 1566

# Row 31: Israel
$ws.Range("B31").Value = 16246
$ws.Range("C31").Value = 38
$ws.Range("D31").Value = 10064
$ws.Range("E31").Value = 5947
$ws.Range("F31").Value = 70
$ws.Range("G31").Value = 3
$ws.Range("H31").Value = 235

# Row 89: Senegal
$ws.Range("F89").Value = 6

# Row 124: now Estado de Palestina - new data
$ws.Range("B124").Value = 362
$ws.Range("C124").Value = 9
$ws.Range("D124").Value = 102
$ws.Range("E124").Value = 258
$ws.Range("F124").Value = 0
$ws.Range("H124").Value = 2

# Row 125: now Venezuela - (unchanged) data
$ws.Range("B125").Value = 357
$ws.Range("C125").Value = 0
$ws.Range("D125").Value = 158
$ws.Range("E125").Value = 189
$ws.Range("F125").Value = 1
$ws.Range("H125").Value = 10

# Row 132: Ruanda
$ws.Range("B132").Value = 261
$ws.Range("C132").Value = 2
$ws.Range("D132").Value = 128
$ws.Range("E132").Value = 133

# --- Update the "last updated" timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 4 de Mayo de 2020 a las 21:33"
